$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 881.5
$ws.Range("J2").Value = 2399.5
$ws.Range("L2").Value = 2399.5
$ws.Range("N2").Value = -2625.5
$ws.Range("H58").Value = 1313.6666
$ws.Range("J58").Value = 1599
$ws.Range("L58").Value = 4797
$ws.Range("N58").Value = -5097
$ws.Range("H88").Value = 6092.316
$ws.Range("I88").Value = 2150
$ws.Range("J88").Value = 6556.1177
$ws.Range("K88").Value = 2150
$ws.Range("L88").Value = 6556.1177
$ws.Range("M88").Value = -1744
$ws.Range("N88").Value = -7368.1177
$ws.Range("H91").Value = 6092.316
$ws.Range("I91").Value = 2150
$ws.Range("J91").Value = 6556.1177
$ws.Range("K91").Value = 2150
$ws.Range("L91").Value = 6556.1177
$ws.Range("M91").Value = -746
$ws.Range("N91").Value = -9364.117699999999
$ws.Range("H129").Value = 8702.083000000001
$ws.Range("I129").Value = 1737.2222
$ws.Range("J129").Value = 29596.666
$ws.Range("K129").Value = 5211.6666
$ws.Range("L129").Value = 88789.99800000001
$ws.Range("M129").Value = -211.6665999999996
$ws.Range("N129").Value = -98789.99800000001
$ws.Range("H132").Value = 2818.3215
$ws.Range("I132").Value = 2889.0386
$ws.Range("K132").Value = 8667.1158
$ws.Range("M132").Value = -6137.1158

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4958.125
$ws.Range("I32").Value = 2869.7778
$ws.Range("K32").Value = 2869.7778
$ws.Range("M32").Value = -2582.7778
$ws.Range("H61").Value = 62501732
$ws.Range("I61").Value = 76924750
$ws.Range("J61").Value = 1999.3334
$ws.Range("K61").Value = 76924750
$ws.Range("L61").Value = 1999.3334
$ws.Range("M61").Value = -76924538
$ws.Range("N61").Value = -2423.3334
$ws.Range("H74").Value = 28574268
$ws.Range("I74").Value = 29414674
$ws.Range("K74").Value = 29414674
$ws.Range("M74").Value = -29413800
$ws.Range("H77").Value = 28574268
$ws.Range("I77").Value = 29414674
$ws.Range("K77").Value = 147073370
$ws.Range("M77").Value = -147069002
$ws.Range("H92").Value = 10000
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -14992
$ws.Range("H103").Value = 100000
$ws.Range("J103").Value = 100000
$ws.Range("L103").Value = 100000
$ws.Range("N103").Value = -102344
$ws.Range("H132").Value = 2502053.5
$ws.Range("I132").Value = 2704715
$ws.Range("K132").Value = 8114145
$ws.Range("M132").Value = -8111615
$ws.Range("H136").Value = 62501732
$ws.Range("I136").Value = 76924750
$ws.Range("J136").Value = 1999.3334
$ws.Range("K136").Value = 230774250
$ws.Range("L136").Value = 5998.0002
$ws.Range("M136").Value = -230771700
$ws.Range("N136").Value = -11098.0002
$ws.Range("H138").Value = 164999
$ws.Range("J138").Value = 164999
$ws.Range("L138").Value = 164999
$ws.Range("N138").Value = -175279

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2247.875
$ws.Range("I86").Value = 2372.9524
$ws.Range("J86").Value = 2009.091
$ws.Range("K86").Value = 2372.9524
$ws.Range("L86").Value = 2009.091
$ws.Range("M86").Value = -1249.9524
$ws.Range("N86").Value = -4255.091
$ws.Range("H89").Value = 2247.875
$ws.Range("I89").Value = 2372.9524
$ws.Range("J89").Value = 2009.091
$ws.Range("K89").Value = 11864.762
$ws.Range("L89").Value = 10045.455
$ws.Range("M89").Value = -6248.762000000001
$ws.Range("N89").Value = -21277.455
$ws.Range("H95").Value = 17162
$ws.Range("J95").Value = 17162
$ws.Range("L95").Value = 17162
$ws.Range("N95").Value = -22654
$ws.Range("H134").Value = 33334492
$ws.Range("I134").Value = 33334492
$ws.Range("K134").Value = 100003476
$ws.Range("M134").Value = -100000941

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3254.2307
$ws.Range("I31").Value = 3595
$ws.Range("K31").Value = 3595
$ws.Range("M31").Value = -3300
$ws.Range("H34").Value = 3254.2307
$ws.Range("I34").Value = 3595
$ws.Range("K34").Value = 3595
$ws.Range("M34").Value = -3393
$ws.Range("H58").Value = 71444940
$ws.Range("I58").Value = 100021520
$ws.Range("K58").Value = 100021520
$ws.Range("M58").Value = -100021317
$ws.Range("H136").Value = 71444940
$ws.Range("I136").Value = 100021520
$ws.Range("K136").Value = 300064560
$ws.Range("M136").Value = -300062010

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 145199
$ws.Range("I114").Value = 500249.5
$ws.Range("J114").Value = 3178.8
$ws.Range("K114").Value = 1500748.5
$ws.Range("L114").Value = 9536.400000000001
$ws.Range("M114").Value = -1497494.5
$ws.Range("N114").Value = -16044.4
$ws.Range("H117").Value = 2579.8333
$ws.Range("I117").Value = 500
$ws.Range("J117").Value = 2995.8
$ws.Range("K117").Value = 1500
$ws.Range("L117").Value = 8987.400000000001
$ws.Range("M117").Value = 1942
$ws.Range("N117").Value = -15871.4
$ws.Range("H129").Value = 3588
$ws.Range("I129").Value = 3134.5
$ws.Range("J129").Value = 4041.5
$ws.Range("K129").Value = 9403.5
$ws.Range("L129").Value = 12124.5
$ws.Range("M129").Value = -4403.5
$ws.Range("N129").Value = -22124.5
$ws.Range("H131").Value = 1927.6364
$ws.Range("J131").Value = 2398.9
$ws.Range("L131").Value = 7196.700000000001
$ws.Range("N131").Value = -17276.7

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 4435
$ws.Range("J92").Value = 4435
$ws.Range("L92").Value = 4435
$ws.Range("N92").Value = -8179

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 180000
$ws.Range("J94").Value = 180000
$ws.Range("L94").Value = 180000
$ws.Range("N94").Value = -181352
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H141").Value = 90830.336
$ws.Range("J141").Value = 90830.336
$ws.Range("L141").Value = 90830.336
$ws.Range("N141").Value = -101190.336

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1971.2307
$ws.Range("I100").Value = 1971.2307
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3942.4614
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3401.4614
$ws.Range("N100").ClearContents()
